# Update column G ("K" - strikeouts replaced with K count) values
# with the regenerated save_data values described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 4
    4  = 3
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 3
    10 = 0
    11 = 1
    12 = 3
    13 = 3
    14 = 2
    15 = 3
    16 = 0
    17 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
